$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = 2023
$ws.Range("H2").Value = 33
$ws.Range("H3").Value = 30
$ws.Range("H4").Value = 65
$ws.Range("H5").Value = 19
$ws.Range("H6").Value = 58
$ws.Range("H7").Value = 29
$ws.Range("H8").Value = 8
$ws.Range("H9").Value = 3
$ws.Range("H10").Value = 24
$ws.Range("H11").Value = 13
$ws.Range("H12").Value = 282

$ws.Range("H2").Select()
